# Daily attendance processing - 2025-11-27 21:45:32
# Normalize the "Recorded By" (column G) lists so that "System" always
# appears first among the comma-separated recorder names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if (-not ($val -is [string])) { continue }

    $parts = $val -split ','
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    # Case-sensitive membership test: "-ccontains"/"-ceq" are unreliable here,
    # so use the .NET string .Equals() instance method (ordinal, case-sensitive).
    $hasExactSystem = $false
    foreach ($p in $trimmed) {
        if ($p.Equals("System")) { $hasExactSystem = $true }
    }

    if ($hasExactSystem) {
        $rest = @()
        foreach ($p in $trimmed) {
            if (-not $p.Equals("System")) { $rest += $p }
        }
        $newParts = @("System") + $rest
        $newVal = $newParts -join ', '

        if (-not $newVal.Equals($val)) {
            $cell.Value = $newVal
        }
    }
}
